# Update cryptocurrency price list (cryptos list refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.380.16"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.36%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.874.54"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.52%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7109"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.78%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.16"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.55%  "

$ws.Range("E7").Value = "  -0.10%  "

$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3115"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.54%  "

$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07778"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.77%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "25.04"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.21%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08463"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.77%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.870.47"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.40%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.236"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.22%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7124"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.84%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.28"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.33%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "29.385.13"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.11%  "

$ws.Range("B17").Value = "Uniswap"
$ws.Range("C17").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.037"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.41%  "

$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008224"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +5.23%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "240.71"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.31%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.25"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.70%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.125.97"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.76%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9996"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.11%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.797"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.97%  "

$ws.Range("E24").Value = "  -0.13%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1599"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.11%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.69"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.49%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.054"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.34%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.49"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.73%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.513"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.15%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.422"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.75%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.319"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.68%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.280"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.87%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05315"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.41%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.937"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.40%  "

$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.177"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.29%  "

$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7496"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.04%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.699"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.67%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01871"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.58%  "

$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.720"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.79%  "

$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.205.36"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.93%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.458"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.62%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "72.98"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.24%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8862"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.05%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "107.95"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.09%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.0000"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.10%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.022.94"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.79%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.817"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.26%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5210"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.22%  "

$ws.Range("E49").Value = "  +8.58%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.413"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.58%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4320"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.80%  "
